$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# --- Step 1: the existing "erreurs" sheet becomes "fréquence_erreurs". This
#     keeps its sheetId/r:id/part (sheet4.xml) so that, exactly like in the
#     target workbook, "fréquence_erreurs" ends up as sheetId 4 / rId4, and a
#     freshly-added sheet right after it becomes sheetId 5 / rId5. ----------
$freqSheet = $wb.Worksheets.Item("erreurs")
$freqSheet.Name = "fréquence_erreurs"

# Wipe whatever is currently on that sheet (the old "erreurs" data) so we can
# lay down the brand-new "fréquence_erreurs" content from scratch.
$freqSheet.Cells.UnMerge()
$freqSheet.Cells.Clear()

# --- Step 2: add a new sheet right after it, named "erreurs", and restore
#     the original "erreurs" content onto it (unchanged vs. the source). ---
$newErr = $wb.Worksheets.Add($null, $freqSheet)
$newErr.Name = "erreurs"

# Row 1
$newErr.Range("B1").Value = "Q_1"
$newErr.Range("D1").Value = "Q_2"
$newErr.Range("F1").Value = "Q_3"
# Row 2
$newErr.Range("B2").Value = "sto"
$newErr.Range("C2").Value = "dv"
$newErr.Range("D2").Value = "rv"
$newErr.Range("E2").Value = "oz"
$newErr.Range("F2").Value = "rm"
# Row 4
$newErr.Range("A4").Value = 12345
$newErr.Range("B4").Value = 1
$newErr.Range("C4").Value = 1
$newErr.Range("D4").Value = 0
$newErr.Range("E4").Value = 1
$newErr.Range("F4").Value = 0
# Row 5
$newErr.Range("A5").Value = 23456
$newErr.Range("B5").Value = 0
$newErr.Range("C5").Value = 1
$newErr.Range("D5").Value = 1
$newErr.Range("E5").Value = 0
$newErr.Range("F5").Value = 1

# Merge first, *then* paste the formatting — pasting formats onto cells that
# get merged afterwards makes Excel synthesize extra "inner border removed"
# style variants; merging first avoids that and keeps everything on the one
# shared "label" style used across the rest of the workbook.
$newErr.Range("B1:C1").Merge()
$newErr.Range("D1:E1").Merge()

# Reuse the existing "label" style (bold / centered / bordered, index 1 in
# the original workbook) by copying its formatting from a cell that already
# carries it, rather than creating a brand-new style entry.
$styleSource = $wb.Worksheets.Item("codes_pondération").Range("A2")
$styleSource.Copy()
$newErr.Range("A1:F2").PasteSpecial($xlPasteFormats)
$newErr.Range("A4:A5").PasteSpecial($xlPasteFormats)

# --- Step 3: populate "fréquence_erreurs" with the new report content. -----
# Row 1 — headers
$freqSheet.Range("C1").Value = "fréquence erreurs (%)"
$freqSheet.Range("D1").Value = "définition"

# Row 2
$freqSheet.Range("A2").Value = "Q_1"
$freqSheet.Range("B2").Value = "dv"
$freqSheet.Range("C2").Value = 100
$freqSheet.Range("D2").Value = "Mauvaise prise en compte de la durée de vie du produit dans le calcul de cycle de vie"

# Row 3
$freqSheet.Range("B3").Value = "sto"
$freqSheet.Range("C3").Value = 50
$freqSheet.Range("D3").Value = "Erreur dans le calcul stoechiométrique du procédé de combustion"

# Row 4
$freqSheet.Range("A4").Value = "Q_2"
$freqSheet.Range("B4").Value = "rv"
$freqSheet.Range("C4").Value = 50
$freqSheet.Range("D4").Value = "Réponse correcte, mais vague"

# Row 5
$freqSheet.Range("B5").Value = "oz"
$freqSheet.Range("C5").Value = 50
$freqSheet.Range("D5").Value = "Confond couche d’ozone et gaz à effet de serre!"

# Row 6
$freqSheet.Range("A6").Value = "Q_3"
$freqSheet.Range("B6").Value = "rm"
$freqSheet.Range("C6").Value = 50
$freqSheet.Range("D6").Value = "réponse manquante (la question n'a pas été faite)"

$freqSheet.Range("A2:A3").Merge()
$freqSheet.Range("A4:A5").Merge()

$styleSource.Copy()
$freqSheet.Range("C1:D1").PasteSpecial($xlPasteFormats)
$freqSheet.Range("A2:B6").PasteSpecial($xlPasteFormats)

$freqSheet.Range("A1").Select()
